# fix: prevent hidden columns from being labeled upon detecting changes
#
# Rows whose FV2410 columns (B:K) are actually identical to their FV2504
# counterparts (M:V) were nonetheless being flagged as "ÄNDERUNG" in column L
# because the comparison logic used to include hidden columns (K/V,
# "Bedingung_FV2410"/"Bedingung_FV2504"). After the fix, such rows no longer
# show the "ÄNDERUNG" label: the L cell is cleared and restyled to the blank
# "no change" look.
#
# Additionally, rows that are the first row of a new "Segmentname" group and
# end up with no detected change get the same gray "group header" highlight
# style that is already used elsewhere in the sheet (e.g. row 9, "Beginn der
# Nachricht").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that start a new Segmentname group AND have no real change left
# -> get the full gray group-header style (copied from row 9, which already
#    carries that exact formatting).
$fullStyleRows = @(13, 17, 23, 27, 34, 40, 71, 99)

# Remaining rows where only the "ÄNDERUNG" flag in column L must be removed.
$lOnlyRows = @(14, 15, 16, 18, 19, 20, 21, 22, 24, 25, 26, 28, 29, 30, 31, 32, 33, 35, 36, 38, 39, 41, 42, 43, 103)

# Copy the whole formatting of the reference "group header" row (row 9) once,
# then paste the format only (not the values) onto each target row.
$ws.Range("A9:V9").Copy()
foreach ($r in $fullStyleRows) {
    $ws.Range("A$r`:V$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

# Copy the blank "no change" formatting used in column L (row 9's L cell),
# then paste the format only onto every affected row's L cell.
$ws.Range("L9").Copy()
foreach ($r in $fullStyleRows) {
    $ws.Range("L$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}
foreach ($r in $lOnlyRows) {
    $ws.Range("L$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

# Finally, clear the "ÄNDERUNG" text itself from every affected row's L cell.
foreach ($r in ($fullStyleRows + $lOnlyRows)) {
    $ws.Range("L$r").ClearContents()
}

Write-Host "Applied hidden-column change-detection fix to rows:"
Write-Host ($fullStyleRows + $lOnlyRows | Sort-Object)
